$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.3096085409252669
$ws1.Range("C2").Value = 0.06097560975609756
$ws1.Range("D2").Value = 0.8928571428571429
$ws1.Range("E2").Value = 0.1141552511415525
$ws1.Range("F2").Value = 0.239463601532567
$ws1.Range("G2").Value = 0.5855855855855856
$ws1.Range("H2").Value = 0.7602327447833066
$ws1.Range("I2").Value = 25
$ws1.Range("J2").Value = 385
$ws1.Range("K2").Value = 149
$ws1.Range("L2").Value = 3

# ---- Sheet: Classification Report ----
$ws2 = $wb.Worksheets.Item("Classification Report")

# Row 2 ("0")
$ws2.Range("B2").Value = 0.9802631578947368
$ws2.Range("C2").Value = 0.2790262172284644
$ws2.Range("D2").Value = 0.434402332361516

# Row 3 ("1")
$ws2.Range("B3").Value = 0.06097560975609756
$ws2.Range("C3").Value = 0.8928571428571429
$ws2.Range("D3").Value = 0.1141552511415525

# Row 4 ("accuracy")
$ws2.Range("B4").Value = 0.3096085409252669
$ws2.Range("C4").Value = 0.3096085409252669
$ws2.Range("D4").Value = 0.3096085409252669
$ws2.Range("E4").Value = 0.3096085409252669

# Row 5 ("macro avg")
$ws2.Range("B5").Value = 0.5206193838254172
$ws2.Range("C5").Value = 0.5859416800428037
$ws2.Range("D5").Value = 0.2742787917515342

# Row 6 ("weighted avg")
$ws2.Range("B6").Value = 0.9344623547846268
$ws2.Range("C6").Value = 0.3096085409252669
$ws2.Range("D6").Value = 0.4184469617669271

# ---- Sheet: Confusion Matrix ----
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 ("Actual 0")
$ws3.Range("B2").Value = 149
$ws3.Range("C2").Value = 385

# Row 3 ("Actual 1")
$ws3.Range("B3").Value = 3
$ws3.Range("C3").Value = 25
